$d = $word.ActiveDocument

# Locate the paragraph that still has the old wording.
$old = "topic_number (*). Would be chapter number.Topic number. E.G. Chapter 1: 1.1, 1.2, 1.2.1 ETC. :Number."
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.TrimEnd([char]13) -eq $old) {
        $target = $para
    }
}

if ($target -eq $null) {
    throw "Could not find target paragraph"
}

$start = $target.Range.Start
$end = $target.Range.End
# Paragraph.Range.End includes the paragraph mark, exclude it so we only
# rewrite the run content and keep the paragraph (and its pPr/numbering) intact.
$contentRange = $d.Range($start, $end - 1)

# Rebuild the sentence as several runs, matching how Word splits text when
# the grammar checker flags a phrase (here "ETC. :Float") with proofErr
# start/end markers, and changes the trailing "Number." to "Float.".
$newXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t xml:space="preserve">topic_number (*). Would be chapter number.Topic number. E.G. Chapter 1: 1.1, 1.2, 1.2.1 </w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:t>ETC. :</w:t></w:r>
<w:r><w:t>Float</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t>.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$contentRange.InsertXML($newXml)
